# Regenerate the "K" column (G) values for the save_data sheet.
# The commit message indicates column G (header "K") values were
# recalculated (it used to mirror "Strike#", now uses "K"); only the
# numeric values in column G change, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G
$newValues = @{
    2  = 1
    3  = 4
    4  = 0
    5  = 2
    6  = 0
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 6
    17 = 4
    18 = 2
    19 = 3
    20 = 8
    21 = 3
    22 = 1
    23 = 0
    24 = 3
    25 = 3
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
